$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("log")

# Append a new log entry as row 95 (the sheet currently has data through row 94)
$ws.Cells.Item(95, 1).Value = 94
$ws.Cells.Item(95, 2).Value = 1029804860
$ws.Cells.Item(95, 3).Value = "Info Pelanggan"
$ws.Cells.Item(95, 4).Value = "2024-06-24 19:16:51.358957"
